# Re-pull data: update column F (dSF) values for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -8
    3  = -2
    5  = 5
    13 = 0
    14 = -4
    18 = 8
    20 = 5
    21 = 7
    27 = -5
    29 = 4
    32 = -5
    34 = 2
    38 = -2
    39 = 3
    43 = 1
    46 = 9
    50 = -2
    51 = -2
    53 = -3
    56 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
